# Slide 1's notes page ("Note from Chuck ...") is translated into Greek,
# keeping the name "Chuck" itself in English in the middle of the sentence.
# In the canonical OOXML the single run becomes five runs (the Greek parts
# tagged lang="el-GR", the name kept lang="en-US", and the "es" ending of
# "σελίδα/ες" flagged err="1"); all keep the note's dk2 scheme-color fill.
#
#   1. el-GR  "Σημείωση από τον "
#   2. en-US  " Chuck"
#   3. el-GR  ". Εάν χρησιμοποιείτε αυτό το υλικό, ... σελίδα/"
#   4. el-GR  "ες"            (err="1")
#   5. el-GR  " αναγνώρισης."

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.NotesPage.Shapes.Item(1)
$tr  = $shp.TextFrame.TextRange

$run1 = "Σημείωση από τον "
$run2 = " Chuck"
$run3 = ". Εάν χρησιμοποιείτε αυτό το υλικό, μπορείτε να αφαιρέσετε το λογότυπο UM και να το αντικαταστήσετε με το δικό σας, αλλά διατηρήστε το λογότυπο CC-BY στην πρώτη σελίδα καθώς την/τις σελίδα/"
$run4 = "ες"
$run5 = " αναγνώρισης."

$tr.Text = $run1 + $run2 + $run3 + $run4 + $run5
